# Applies the automated sales update described by the commit
# "Actualización automática 2025-09-10 16:05:09".
#
# Changes:
#  - Sheet "VENTAS POR GRUPO": M3 (PORCELANATO) 0 -> 43.1
#  - Sheet "VENTAS POR GRUPO": M10 (progress label) "1 de 8" -> "2 de 8"
#  - Sheet "VENTA MENSUAL":   F3 (septiembre)     0 -> 43.1
#  - Sheet "VENTA MENSUAL":   F10 (total septiembre) 5372.02 -> 5415.120000000001

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M3").Value = 43.1
$wsGrupo.Range("M10").Value = "2 de 8"

$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 43.1
$wsMensual.Range("F10").Value = 5415.120000000001
